$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 10.5
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 60
$ws.Range("M6").Value = 109
$ws.Range("N6").Value = -284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 5865.3335
$ws.Range("I4").Value = 4397.5
$ws.Range("J4").Value = 6599.25
$ws.Range("K4").Value = 4397.5
$ws.Range("L4").Value = 6599.25
$ws.Range("M4").Value = -4281.5

$ws.Range("H63").Value = 27649
$ws.Range("I63").Value = 41415
$ws.Range("J63").Value = 7000
$ws.Range("K63").Value = 41415
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = -40729
$ws.Range("N63").Value = -8372

$ws.Range("H66").Value = 27649
$ws.Range("I66").Value = 41415
$ws.Range("J66").Value = 7000
$ws.Range("K66").Value = 207075
$ws.Range("L66").Value = 35000
$ws.Range("M66").Value = -203643
$ws.Range("N66").Value = -41864

$ws.Range("H74").Value = 1761.7894
$ws.Range("I74").Value = 1865.8182
$ws.Range("J74").Value = 1618.75
$ws.Range("K74").Value = 1865.8182
$ws.Range("L74").Value = 1618.75
$ws.Range("M74").Value = -991.8181999999999

$ws.Range("H77").Value = 1761.7894
$ws.Range("I77").Value = 1865.8182
$ws.Range("J77").Value = 1618.75
$ws.Range("K77").Value = 9329.091
$ws.Range("L77").Value = 8093.75
$ws.Range("M77").Value = -4961.091

$ws.Range("H97").Value = 1632.6666
$ws.Range("I97").Value = 1779.4
$ws.Range("J97").Value = 899
$ws.Range("K97").Value = 1779.4
$ws.Range("L97").Value = 899
$ws.Range("M97").Value = -1283.4
$ws.Range("N97").Value = -1891

$ws.Range("H123").Value = 43476
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 43476
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 43476
$ws.Range("N123").Value = -53276

$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 3099.5
$ws.Range("I132").Value = 3099.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9298.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6768.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 652.75
$ws.Range("I22").Value = 511.7143
$ws.Range("J22").Value = 1640
$ws.Range("K22").Value = 511.7143
$ws.Range("L22").Value = 1640
$ws.Range("M22").Value = -338.7143

$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -877

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4384

$ws.Range("H94").Value = 1464.6666
$ws.Range("I94").Value = 1314.1666
$ws.Range("J94").Value = 2066.6667
$ws.Range("K94").Value = 1314.1666
$ws.Range("L94").Value = 2066.6667
$ws.Range("M94").Value = -863.1666
$ws.Range("N94").Value = -2968.6667

$ws.Range("H99").Value = 2150.9167
$ws.Range("I99").Value = 1868
$ws.Range("J99").Value = 2999.6667
$ws.Range("K99").Value = 1868
$ws.Range("L99").Value = 2999.6667
$ws.Range("M99").Value = -370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 980
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -550
$ws.Range("N22").Value = -1700

$ws.Range("H54").Value = 38061.5
$ws.Range("I54").Value = 37388.668
$ws.Range("J54").Value = 40080
$ws.Range("K54").Value = 37388.668
$ws.Range("L54").Value = 40080
$ws.Range("M54").Value = -36730.668
$ws.Range("N54").Value = -41396

$ws.Range("H58").Value = 1919.6842
$ws.Range("I58").Value = 1733.8235
$ws.Range("J58").Value = 3499.5
$ws.Range("K58").Value = 1733.8235
$ws.Range("L58").Value = 3499.5
$ws.Range("M58").Value = -1530.8235

$ws.Range("H94").Value = 74959.31
$ws.Range("I94").Value = 127780.664
$ws.Range("J94").Value = 7046.143
$ws.Range("K94").Value = 127780.664
$ws.Range("L94").Value = 7046.143
$ws.Range("M94").Value = -127329.664
$ws.Range("N94").Value = -7948.143

$ws.Range("N114").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0

$ws.Range("H136").Value = 1919.6842
$ws.Range("I136").Value = 1733.8235
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 5201.470499999999
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -2651.470499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11000062
$ws.Range("I4").Value = 11000062
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 33000186
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -33000074

$ws.Range("H12").Value = 9557.6
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 10614
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 31842
$ws.Range("M12").Value = 23
$ws.Range("N12").Value = -32188

$ws.Range("H23").Value = 75000070
$ws.Range("I23").Value = 100000080
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = 300000240
$ws.Range("L23").Value = 114
$ws.Range("M23").Value = -300000005

$ws.Range("H60").Value = 654.3333
$ws.Range("I60").Value = 406.4
$ws.Range("J60").Value = 1894
$ws.Range("K60").Value = 1219.2
$ws.Range("L60").Value = 5682
$ws.Range("M60").Value = -968.1999999999998
$ws.Range("N60").Value = -6184

$ws.Range("H64").Value = 18895
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 18895
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 56685
$ws.Range("N64").Value = -57225

$ws.Range("H67").Value = 18895
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 18895
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 56685
$ws.Range("N67").Value = -58557

$ws.Range("H103").Value = 2467
$ws.Range("I103").Value = 1425
$ws.Range("J103").Value = 3509
$ws.Range("K103").Value = 4275
$ws.Range("L103").Value = 10527
$ws.Range("M103").Value = -3396
$ws.Range("N103").Value = -12285

$ws.Range("H117").Value = 17007.666
$ws.Range("I117").Value = 409.2
$ws.Range("J117").Value = 100000
$ws.Range("K117").Value = 1227.6
$ws.Range("L117").Value = 300000
$ws.Range("M117").Value = 2214.4

$ws.Range("H121").Value = 8252.727999999999
$ws.Range("I121").Value = 25577.6
$ws.Range("J121").Value = 3157.1765
$ws.Range("K121").Value = 76732.79999999999
$ws.Range("L121").Value = 9471.529500000001
$ws.Range("M121").Value = -75422.79999999999
$ws.Range("N121").Value = -12091.5295

$ws.Range("H140").Value = 5594.9443
$ws.Range("I140").Value = 1977.8462
$ws.Range("J140").Value = 14999.4
$ws.Range("K140").Value = 5933.5386
$ws.Range("L140").Value = 44998.2
$ws.Range("M140").Value = -753.5385999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 161.8
$ws.Range("I2").Value = 178.11111
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 178.11111
$ws.Range("L2").Value = 15
$ws.Range("M2").Value = -65.11111
$ws.Range("N2").Value = -241

$ws.Range("H70").Value = 8708
$ws.Range("I70").Value = 6960.364
$ws.Range("J70").Value = 9989.6
$ws.Range("K70").Value = 6960.364
$ws.Range("L70").Value = 9989.6
$ws.Range("M70").Value = -6690.364
$ws.Range("N70").Value = -10529.6

$ws.Range("H73").Value = 8708
$ws.Range("I73").Value = 6960.364
$ws.Range("J73").Value = 9989.6
$ws.Range("K73").Value = 6960.364
$ws.Range("L73").Value = 9989.6
$ws.Range("M73").Value = -6024.364
$ws.Range("N73").Value = -11861.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3518.3333
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 3518.3333
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3518.3333
$ws.Range("N16").Value = -3858.3333

$ws.Range("H22").Value = 1620
$ws.Range("I22").Value = 2473.75
$ws.Range("J22").Value = 1132.1428
$ws.Range("K22").Value = 2473.75
$ws.Range("L22").Value = 1132.1428
$ws.Range("M22").Value = -2178.75
$ws.Range("N22").Value = -1722.1428

$ws.Range("H27").Value = 1620
$ws.Range("I27").Value = 2473.75
$ws.Range("J27").Value = 1132.1428
$ws.Range("K27").Value = 2473.75
$ws.Range("L27").Value = 1132.1428
$ws.Range("M27").Value = -2366.75
$ws.Range("N27").Value = -1346.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20544.834
$ws.Range("I62").Value = 26263.166
$ws.Range("J62").Value = 14826.5
$ws.Range("K62").Value = 26263.166
$ws.Range("L62").Value = 14826.5
$ws.Range("M62").Value = -25639.166
$ws.Range("N62").Value = -16074.5

$ws.Range("H65").Value = 20544.834
$ws.Range("I65").Value = 26263.166
$ws.Range("J65").Value = 14826.5
$ws.Range("K65").Value = 131315.83
$ws.Range("L65").Value = 74132.5
$ws.Range("M65").Value = -128195.83
$ws.Range("N65").Value = -80372.5

$ws.Range("H96").Value = 1642.2858
$ws.Range("I96").Value = 1498
$ws.Range("J96").Value = 1700
$ws.Range("K96").Value = 1498
$ws.Range("L96").Value = 1700
$ws.Range("M96").Value = -125
$ws.Range("N96").Value = -4446

$ws.Range("H113").Value = 963.5714
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 581.6667
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 1745.0001
$ws.Range("M113").Value = -1580
$ws.Range("N113").Value = -6085.0001

$ws.Range("H140").Value = 84997.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 84997.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 84997.5
$ws.Range("N140").Value = -95357.5

$ws.Range("H141").Value = 179966.67
$ws.Range("I141").Value = 89900
$ws.Range("J141").Value = 225000
$ws.Range("K141").Value = 89900
$ws.Range("L141").Value = 225000
$ws.Range("M141").Value = -84720
$ws.Range("N141").Value = -235360
